$wb = $excel.ActiveWorkbook

# --- "Weekly Quantity" sheet ---
$ws1 = $wb.Worksheets.Item("Weekly Quantity")

# Delete rows 14-17 (they are being merged/removed; row 13 now carries the
# date+qty that used to live in row 17).
$ws1.Rows("14:17").Delete()

# Row 13: date moves to the last week's date (45137.99999999999); quantity
# becomes 540 (the value that had been on the now-deleted row 17).
$ws1.Range("A13").Value = 45137.99999999999
$ws1.Range("B13").Value = 540

# --- "Monthly Trend" sheet ---
$ws2 = $wb.Worksheets.Item("Monthly Trend")
$ws2.Range("B6").Value = 10
$ws2.Range("B7").Value = 540

Write-Output "edit applied"
